$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2014")

# Row 13: the "157" weight was actually an uncalibrated scale reading -> store as text note
$ws.Range("K13").Value = "157 scale no calibrated "

# Row 15: the "218" weight was actually an uncalibrated scale reading -> store as text note
$ws.Range("K15").Value = "218 Scale not Calibrated "

# Row 22: remove the legend note about gray cells / calibration
$ws.Range("I22").Value = ""

# Row 28: split the combined "species_sex" header into separate "species" and "sex" columns,
# and move the "wgt" header over to column F
$ws.Range("D28").Value = "species "
$ws.Range("E28").Value = "sex "
$ws.Range("F28").Value = "wgt"

# Rows 29-33: split the combined species_sex codes (e.g. "DM_F") into separate species and
# sex columns, and shift the weight value from column E to column F
$ws.Range("D29").Value = "DM"
$ws.Range("E29").Value = "F"
$ws.Range("F29").Value = 37

$ws.Range("D30").Value = "DM"
$ws.Range("E30").Value = "M"
$ws.Range("F30").Value = ""

$ws.Range("D31").Value = "DM"
$ws.Range("E31").Value = "F "
$ws.Range("F31").Value = 48

$ws.Range("D32").Value = "DO"
$ws.Range("E32").Value = "M"
$ws.Range("F32").Value = 52

$ws.Range("D33").Value = "OL"
$ws.Range("E33").Value = "M"
$ws.Range("F33").Value = 35
